# Insert a new row at position 43 (pushes existing rows 43-53 down to 44-54,
# carrying their values/styles with them, matching Excel's native Insert
# behaviour), then populate the newly inserted row with the
# "planned pharmacological substance dose" (GMHO:0000200) record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = "GMHO:0000200"
$ws.Range("B43").Value = "planned pharmacological substance dose"
$ws.Range("C43").Value = "A plan specification about the dose of pharmacological substance."
$ws.Range("D43").Value = "plan specification"
$ws.Range("P43").Value = "LSR 1"
$ws.Range("Q43").Value = "Intervention content and delivery"
$ws.Range("S43").Value = "Proposed"
$ws.Range("V43").Value = "PS"
